$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Fill in previously-empty "Resolucion" cells with "Ok" ---
# Row 14: Envido + Envido + Real + F | Quiero | (empty) -> Ok
$t.Cell(14, 3).Range.Text = "Ok"
# Row 18: Real + Falta | Quiero | (empty) -> Ok
$t.Cell(18, 3).Range.Text = "Ok"
# Row 20: Falta | Quiero | (empty) -> Ok
$t.Cell(20, 3).Range.Text = "Ok"
# Row 29: Flor | Aceptar | (empty) -> Ok
$t.Cell(29, 3).Range.Text = "Ok"

# Row 30: Envido + Flor por a... | (empty) | Ok -> Aceptar in 2nd cell
$t.Cell(30, 2).Range.Text = "Aceptar"

# --- Insert 4 new rows right before row 30 ("Envido + Flor por a...") ---
$targetRow = $t.Rows.Item(30)
$t.Rows.Add($targetRow) | Out-Null
$t.Rows.Add($targetRow) | Out-Null
$t.Rows.Add($targetRow) | Out-Null
$t.Rows.Add($targetRow) | Out-Null

# Fill the 4 new rows (now at indices 30-33) with their content.
# Cells that must stay empty are first given a single placeholder
# character and then that character is deleted, which collapses the
# paragraph back down to a clean, run-less <w:p/> (matching how Word
# represents an untouched empty cell) instead of leaving a stray empty
# <w:r/> behind (which is what a freshly-added row's cell looks like by
# default).
$t.Cell(30, 1).Range.Text = "Flor + flor"
$t.Cell(30, 2).Range.Text = "Quiero"
$t.Cell(30, 3).Range.Text = "X"
$t.Cell(30, 3).Range.Characters.Item(1).Delete()

$t.Cell(31, 1).Range.Text = "Flor + flor"
$t.Cell(31, 2).Range.Text = "Flor me achico"
$t.Cell(31, 3).Range.Text = "Ok"

$t.Cell(32, 1).Range.Text = "Flor + contra flor"
$t.Cell(32, 2).Range.Text = "Quiero"
$t.Cell(32, 3).Range.Text = "X"
$t.Cell(32, 3).Range.Characters.Item(1).Delete()

$t.Cell(33, 1).Range.Text = "Flor + contra flor"
$t.Cell(33, 2).Range.Text = "No"
$t.Cell(33, 3).Range.Text = "X"
$t.Cell(33, 3).Range.Characters.Item(1).Delete()
